$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D17:D28").Value = "Implementation"
$ws.Range("C31").Value = "Module Entity/Content (From DB Tables)"

$ws.Range("A7").Select()
$ws.Range("C31").Select()
